# Updated symbol list on Fri Dec 16 05:12:39 UTC 2022 with GitHub Actions
# Updates the "Price" (column D) for coins whose price changed, and bumps
# the "Hora" (column G) value from 4 to 5 for every data row (2-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" values (column D) for the rows that changed.
$dValues = @{
    2  = "261.11"
    3  = "22.94"
    4  = "6.208"
    5  = "0.06141"
    6  = "6.735"
    7  = "3.459"
    8  = "1.365"
    9  = "0.7982"
    10 = "0.1574"
    11 = "0.08130"
    12 = "0.03496"
    13 = "0.03085"
    14 = "0.09323"
    15 = "3.879"
    16 = "0.001689"
    17 = "0.04782"
    18 = "0.0006143"
    19 = "0.006208"
    21 = "0.004066"
    22 = "0.0001501"
    23 = "3.692"
    24 = "2.212"
    25 = "0.3362"
    27 = "0.0003204"
    40 = "0.04607"
    41 = "0.007083"
    42 = "0.1120"
    43 = "0.003161"
    46 = "0.00005935"
    48 = "0.7003"
    49 = "0.08931"
}

foreach ($row in $dValues.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $dValues[$row]
}

# Bump "Hora" (column G) from 4 to 5 for every data row.
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $cell.NumberFormat = "@"
    $cell.Value = "5"
}
